{"js": "// Remove ALL horizontal-line paragraphs (the zero-width VML <v:rect>\n// \"hr\" pictures Word inserts as their own empty paragraph) from the\n// document body. Such a paragraph carries no real text \u2014 only the\n// drawn rule \u2014 so it reads back with empty text; confirm it really is\n// a horizontal-rule paragraph by checking its OOXML for the `o:hr`\n// marker before deleting it (guards against deleting an unrelated\n// blank paragraph that happens to have no text).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst candidates = paragraphs.items.filter((p) => p.text.length === 0);\n\nconst toDelete = [];\nfor (const p of candidates) {\n  const ooxml = p.getOoxml();\n  await context.sync();\n  if (/\\bo:hr=\"t\"/.test(ooxml.value)) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Remove ALL horizontal-line paragraphs (the zero-width VML <v:rect>\n# \"hr\" pictures Word inserts as their own empty paragraph) from the\n# document. Such a paragraph holds no real text -- just the drawn rule\n# -- so its Range.Text is only the trailing paragraph mark (length 1).\n# Confirm it is really a horizontal rule (and not some unrelated blank\n# paragraph) by checking the paragraph's underlying OOXML for the\n# `o:hr` marker before deleting it. Walk back-to-front so deleting one\n# paragraph doesn't disturb the indices of paragraphs not yet visited.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Length -le 1 -and $p.Range.WordOpenXML -like '*o:hr=*') {\n        $p.Range.Delete()\n    }\n}\n"}
